# In the "Rules" sheet, row 10 (rule R30) the "Integer min" / From value
# (cell C10) is restored from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
